# PastAusElectionPendulums.xlsx - "Update Pendulum with 2019 data"
#
# The 2019 sheet's formulas (K/L/M/... columns) all reference a summary
# row at $153 ($C$153, $E$153, $J$153) that was missing from the sheet,
# so every dependent calculation was resolving against blank cells. This
# change fills in that missing aggregate row (the national swing figures)
# and extends the last electorate row (152, "Wright") with its own
# J/K/L/M pendulum calculations - matching the pattern already used by
# every other election-year sheet in the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 2019 sheet: add the missing J/K/L/M formulas on the last electorate
# row (152) and the new aggregate row (153).
# ---------------------------------------------------------------------
$ws2019 = $wb.Worksheets.Item("2019")

$ws2019.Range("J152").Formula = '=ROUND($E152+IF($I152="ALP",-1,IF(OR($I152="LIB",$I152="NAT",$I152="LIB/NAT"),1,0))-IF($H152="ALP",-1,IF(OR($H152="LIB",$H152="NAT",$H152="LIB/NAT"),1,0)),2)'
$ws2019.Range("K152").Formula = '=ROUND($J152+($E$153-$J$153),2)'
$ws2019.Range("L152").Formula = '=IF($P$7="Pre-election",ROUND($K152+($P$8-$E$153),2),ROUND($C152+($P$8-$C$153),2))'
$ws2019.Range("M152").Formula = '=_xlfn.NORM.DIST(0,$L152,3.3,TRUE)'

$ws2019.Range("C153").Value = 1.53
$ws2019.Range("E153").Value = 0.36
$ws2019.Range("J153").Formula = '=ROUND(SUMPRODUCT($B$2:$B152,J$2:J152)/SUM($B$2:$B152),2)'

# ---------------------------------------------------------------------
# Leftover view-state (selection / scroll position) from the author
# browsing each year's sheet while updating the data. These don't
# change any values - only what's selected/visible when the sheet is
# next opened.
# ---------------------------------------------------------------------
$ws1993 = $wb.Worksheets.Item("1993")
$ws1993.Activate()
$ws1993.Range("C1:I1").Select()

$ws1996 = $wb.Worksheets.Item("1996")
$ws1996.Activate()
$ws1996.Range("C2:I149").Select()

$ws1998 = $wb.Worksheets.Item("1998")
$ws1998.Activate()
$ws1998.Range("C2:I148").Select()

$ws2001 = $wb.Worksheets.Item("2001")
$ws2001.Activate()
$ws2001.Range("C153").Select()

$ws2004 = $wb.Worksheets.Item("2004")
$ws2004.Activate()
$ws2004.Range("C2:I151").Select()

$ws2007 = $wb.Worksheets.Item("2007")
$ws2007.Activate()
$ws2007.Range("C2:I151").Select()

$ws2010 = $wb.Worksheets.Item("2010")
$ws2010.Activate()
$ws2010.Range("C2:I151").Select()

$ws2013 = $wb.Worksheets.Item("2013")
$ws2013.Activate()
$ws2013.Range("C2:I151").Select()

$ws2016 = $wb.Worksheets.Item("2016")
$ws2016.Activate()
$ws2016.Range("J152").Select()

# Restore the originally-active sheet so the workbook reopens on
# "Summary" as before.
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()
